$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New interleaved data for rows 17-28: YOVANIS (CC 78756892) and OSCAR
# (CC 1032415619) alternate per period, in ascending period order
# (1702, 1704, 1705, 1706, 1707, 1708), replacing the previous layout
# where each worker's six periods were grouped together in descending
# order.

$rows = @(
  @{ Row = 17; C = "78756892";   D = "YOVANIS ANTONIO VIGA OSORIO";      E = "1702"; F = 29509;  G = 781242 },
  @{ Row = 18; C = "1032415619"; D = "OSCAR JAVIER RONCANCIO VALBUENA";  E = "1702"; F = 140000; G = 3500000 },
  @{ Row = 19; C = "78756892";   D = "YOVANIS ANTONIO VIGA OSORIO";      E = "1704"; F = 29509;  G = 781242 },
  @{ Row = 20; C = "1032415619"; D = "OSCAR JAVIER RONCANCIO VALBUENA";  E = "1704"; F = 140000; G = 3500000 },
  @{ Row = 21; C = "78756892";   D = "YOVANIS ANTONIO VIGA OSORIO";      E = "1705"; F = 29509;  G = 781242 },
  @{ Row = 22; C = "1032415619"; D = "OSCAR JAVIER RONCANCIO VALBUENA";  E = "1705"; F = 140000; G = 3500000 },
  @{ Row = 23; C = "78756892";   D = "YOVANIS ANTONIO VIGA OSORIO";      E = "1706"; F = 29509;  G = 781242 },
  @{ Row = 24; C = "1032415619"; D = "OSCAR JAVIER RONCANCIO VALBUENA";  E = "1706"; F = 140000; G = 3500000 },
  @{ Row = 25; C = "78756892";   D = "YOVANIS ANTONIO VIGA OSORIO";      E = "1707"; F = 29509;  G = 781242 },
  @{ Row = 26; C = "1032415619"; D = "OSCAR JAVIER RONCANCIO VALBUENA";  E = "1707"; F = 140000; G = 3500000 },
  @{ Row = 27; C = "78756892";   D = "YOVANIS ANTONIO VIGA OSORIO";      E = "1708"; F = 29509;  G = 781242 },
  @{ Row = 28; C = "1032415619"; D = "OSCAR JAVIER RONCANCIO VALBUENA";  E = "1708"; F = 140000; G = 3500000 }
)

foreach ($r in $rows) {
  $i = $r.Row
  $ws.Range("C$i").Value = $r.C
  $ws.Range("D$i").Value = $r.D
  $ws.Range("E$i").Value = $r.E
  $ws.Range("F$i").Value = $r.F
  $ws.Range("G$i").Value = $r.G
}
